# This script re-shuffles the data rows (2-24) of the sheet across the
# columns D, H, I, J, K, L, M, P (Fecha, Variedad, Calidad, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg).
# The other columns (A, B, C, E, F, G, N, O, Q, R) are identical for every
# row in this workbook, so they do not need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values currently found at the
# source row, before any edits, should end up at the destination row).
$rowMap = @{
    2  = 19
    3  = 14
    4  = 20
    5  = 6
    6  = 12
    7  = 23
    8  = 22
    9  = 17
    10 = 8
    11 = 11
    12 = 7
    13 = 10
    14 = 13
    15 = 15
    16 = 4
    17 = 9
    18 = 24
    19 = 16
    20 = 3
    21 = 18
    22 = 5
    23 = 21
    24 = 2
}

# Columns that get shuffled between rows.
$cols = @("D", "H", "I", "J", "K", "L", "M", "P")

# Snapshot the original values for every row/column involved before writing
# anything, since several rows swap values with each other.
# NOTE: reading via the bare ".Value" property is unreliable in this
# COM-interop runtime (it can yield the property descriptor instead of the
# cell contents), so ".Value2" is used for reads instead.
$original = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $original.ContainsKey($srcRow)) {
        $rowValues = @{}
        foreach ($col in $cols) {
            $rowValues[$col] = $ws.Range("$col$srcRow").Value2
        }
        $original[$srcRow] = $rowValues
    }
}

# Apply the permutation.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowValues = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $rowValues[$col]
    }
}
